# Update the cryptocurrency price/volume snapshot data (rows 2-51 of Sheet1).
# Numeric-looking "Price" strings are forced to Text format before assignment
# so Excel keeps them as literal strings (e.g. "1.000", "2.260") instead of
# auto-converting them to numbers and dropping significant digits.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '28.663.04'
$ws.Range('E2').Value = '  +4.31%  '
$ws.Range('D3').Value = '1.874.14'
$ws.Range('E3').Value = '  +2.51%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '338.41'
$ws.Range('E5').Value = '  +2.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9991'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4706'
$ws.Range('E7').Value = '  +3.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4005'
$ws.Range('E8').Value = '  +5.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.73'
$ws.Range('E9').Value = '  +2.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08055'
$ws.Range('E10').Value = '  +2.28%  '
$ws.Range('E11').Value = '  +3.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.18'
$ws.Range('E12').Value = '  +5.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.070'
$ws.Range('E13').Value = '  +3.38%  '
$ws.Range('D14').Value = '1.865.95'
$ws.Range('E14').Value = '  +2.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.297'
$ws.Range('E15').Value = '  +3.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.68'
$ws.Range('E16').Value = '  +1.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06614'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('E20').Value = '  +2.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('D22').Value = '28.674.55'
$ws.Range('E22').Value = '  +4.47%  '
$ws.Range('E23').Value = '  +3.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.06'
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.260'
$ws.Range('E25').Value = '  -2.07%  '
$ws.Range('D26').Value = '2.085.15'
$ws.Range('E26').Value = '  +2.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.59'
$ws.Range('E27').Value = '  +2.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.84'
$ws.Range('E28').Value = '  +2.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.133'
$ws.Range('E29').Value = '  +3.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.522'
$ws.Range('E30').Value = '  +5.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '120.15'
$ws.Range('E31').Value = '  +1.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9872'
$ws.Range('E32').Value = '  +4.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09550'
$ws.Range('E33').Value = '  +2.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.669'
$ws.Range('E34').Value = '  +2.48%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.389'
$ws.Range('E35').Value = '  +4.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.386'
$ws.Range('E36').Value = '  +2.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06201'
$ws.Range('E37').Value = '  +4.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02267'
$ws.Range('E38').Value = '  +4.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.488'
$ws.Range('E39').Value = '  +5.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.186'
$ws.Range('E40').Value = '  +2.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5971'
$ws.Range('E41').Value = '  +3.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9991'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1890'
$ws.Range('E43').Value = '  +3.27%  '
$ws.Range('E44').Value = '  +3.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.261'
$ws.Range('E45').Value = '  -1.54%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.28'
$ws.Range('E46').Value = '  +2.40%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5588'
$ws.Range('E47').Value = '  +2.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.07397'
$ws.Range('E48').Value = '  +11.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.965'
$ws.Range('E49').Value = '  +5.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.079'
$ws.Range('E50').Value = '  +12.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '112.59'
$ws.Range('E51').Value = '  +2.04%  '
